$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44497
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 13000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 14000
$ws.Range("P2").Value = 560

# Row 3
$ws.Range("D3").Value = 44412
$ws.Range("K3").Value = 25000
$ws.Range("L3").Value = 27000
$ws.Range("M3").Value = 26000
$ws.Range("P3").Value = 1040

# Row 4
$ws.Range("D4").Value = 44503
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 11000
$ws.Range("L4").Value = 13000
$ws.Range("M4").Value = 12000
$ws.Range("P4").Value = 480

# Row 5
$ws.Range("D5").Value = 44448
$ws.Range("J5").Value = 400

# Row 6
$ws.Range("D6").Value = 44371
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 28000
$ws.Range("L6").Value = 30000
$ws.Range("M6").Value = 29000
$ws.Range("P6").Value = 1160

# Row 7
$ws.Range("D7").Value = 44364
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 28000
$ws.Range("L7").Value = 30000
$ws.Range("M7").Value = 29000
$ws.Range("P7").Value = 1160

# Row 8
$ws.Range("D8").Value = 44419
$ws.Range("J8").Value = 600
$ws.Range("K8").Value = 27000
$ws.Range("L8").Value = 29000
$ws.Range("M8").Value = 28000
$ws.Range("P8").Value = 1120

# Row 9
$ws.Range("D9").Value = 44447
$ws.Range("J9").Value = 600
$ws.Range("K9").Value = 28000
$ws.Range("L9").Value = 30000
$ws.Range("M9").Value = 29000
$ws.Range("P9").Value = 1160

# Row 10
$ws.Range("D10").Value = 44468
$ws.Range("K10").Value = 23000
$ws.Range("L10").Value = 25000
$ws.Range("M10").Value = 24000
$ws.Range("P10").Value = 960

# Row 11
$ws.Range("D11").Value = 44476
$ws.Range("J11").Value = 500
$ws.Range("K11").Value = 23000
$ws.Range("L11").Value = 24000
$ws.Range("M11").Value = 23500
$ws.Range("P11").Value = 940

# Row 12
$ws.Range("D12").Value = 44427
$ws.Range("J12").Value = 300

# Row 13
$ws.Range("D13").Value = 44434
$ws.Range("K13").Value = 28000
$ws.Range("L13").Value = 30000
$ws.Range("M13").Value = 29000
$ws.Range("P13").Value = 1160

# Row 14
$ws.Range("D14").Value = 44391
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 26000
$ws.Range("L14").Value = 28000
$ws.Range("M14").Value = 27000
$ws.Range("P14").Value = 1080

# Row 15
$ws.Range("D15").Value = 44405
$ws.Range("K15").Value = 26000
$ws.Range("L15").Value = 28000
$ws.Range("M15").Value = 27000
$ws.Range("P15").Value = 1080

# Row 16
$ws.Range("D16").Value = 44377
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 26000
$ws.Range("L16").Value = 28000
$ws.Range("M16").Value = 27000
$ws.Range("P16").Value = 1080

# Row 17
$ws.Range("D17").Value = 44441
$ws.Range("J17").Value = 700
$ws.Range("K17").Value = 28000
$ws.Range("L17").Value = 30000
$ws.Range("M17").Value = 29000
$ws.Range("P17").Value = 1160

# Row 18
$ws.Range("D18").Value = 44384
$ws.Range("K18").Value = 26000
$ws.Range("M18").Value = 27000
$ws.Range("P18").Value = 1080

# Row 20
$ws.Range("D20").Value = 44435
$ws.Range("J20").Value = 900
$ws.Range("K20").Value = 28000
$ws.Range("L20").Value = 30000
$ws.Range("M20").Value = 29000
$ws.Range("P20").Value = 1160

# Row 21
$ws.Range("D21").Value = 44406
$ws.Range("J21").Value = 600

# Row 22
$ws.Range("D22").Value = 44350
$ws.Range("J22").Value = 700

# Row 23
$ws.Range("D23").Value = 44398
$ws.Range("J23").Value = 500
$ws.Range("K23").Value = 26000
$ws.Range("L23").Value = 28000
$ws.Range("M23").Value = 27000
$ws.Range("P23").Value = 1080

# Row 24
$ws.Range("D24").Value = 44392
$ws.Range("J24").Value = 100

# Row 25
$ws.Range("D25").Value = 44433
$ws.Range("J25").Value = 400
$ws.Range("K25").Value = 28000
$ws.Range("L25").Value = 30000
$ws.Range("M25").Value = 29000
$ws.Range("P25").Value = 1160

# Row 26
$ws.Range("D26").Value = 44363
$ws.Range("J26").Value = 240

# Row 27
$ws.Range("D27").Value = 44461
$ws.Range("L27").Value = 25000
$ws.Range("M27").Value = 24000
$ws.Range("P27").Value = 960

# Row 28
$ws.Range("D28").Value = 44475
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 22000
$ws.Range("L28").Value = 24000
$ws.Range("M28").Value = 23000
$ws.Range("P28").Value = 920

# Row 30
$ws.Range("D30").Value = 44482
$ws.Range("J30").Value = 500
$ws.Range("K30").Value = 18000
$ws.Range("L30").Value = 20000
$ws.Range("M30").Value = 19000
$ws.Range("P30").Value = 760

# Row 31
$ws.Range("D31").Value = 44454
$ws.Range("J31").Value = 1000
$ws.Range("K31").Value = 28000
$ws.Range("L31").Value = 30000
$ws.Range("M31").Value = 29000
$ws.Range("P31").Value = 1160

# Row 32
$ws.Range("D32").Value = 44420
$ws.Range("J32").Value = 700
$ws.Range("K32").Value = 27000
$ws.Range("L32").Value = 29000
$ws.Range("M32").Value = 28000
$ws.Range("P32").Value = 1120

# Row 33
$ws.Range("D33").Value = 44343
$ws.Range("K33").Value = 26000
$ws.Range("L33").Value = 28000
$ws.Range("M33").Value = 27000
$ws.Range("P33").Value = 1080

# Row 35
$ws.Range("D35").Value = 44357
$ws.Range("J35").Value = 340
$ws.Range("K35").Value = 28000
$ws.Range("L35").Value = 30000
$ws.Range("M35").Value = 29000
$ws.Range("P35").Value = 1160

# Row 36
$ws.Range("D36").Value = 44455
$ws.Range("J36").Value = 800

# Row 37
$ws.Range("D37").Value = 44490
$ws.Range("J37").Value = 500
$ws.Range("K37").Value = 16000
$ws.Range("L37").Value = 18000
$ws.Range("M37").Value = 17000
$ws.Range("P37").Value = 680

# Row 38
$ws.Range("D38").Value = 44462
$ws.Range("J38").Value = 400
$ws.Range("K38").Value = 22000
$ws.Range("L38").Value = 23000
$ws.Range("M38").Value = 22500
$ws.Range("P38").Value = 900

# Row 39
$ws.Range("D39").Value = 44483
$ws.Range("J39").Value = 300
$ws.Range("K39").Value = 18000
$ws.Range("L39").Value = 20000
$ws.Range("M39").Value = 19000
$ws.Range("P39").Value = 760

# Row 41
$ws.Range("D41").Value = 44504
$ws.Range("J41").Value = 600
$ws.Range("K41").Value = 11000
$ws.Range("L41").Value = 13000
$ws.Range("M41").Value = 12000
$ws.Range("P41").Value = 480

# Row 42
$ws.Range("D42").Value = 44370
$ws.Range("J42").Value = 400
$ws.Range("K42").Value = 27000
$ws.Range("L42").Value = 28000
$ws.Range("M42").Value = 27500
$ws.Range("P42").Value = 1100

# Row 43
$ws.Range("D43").Value = 44385
$ws.Range("J43").Value = 500
$ws.Range("K43").Value = 26000
$ws.Range("L43").Value = 28000
$ws.Range("M43").Value = 27000
$ws.Range("P43").Value = 1080

# Row 44
$ws.Range("D44").Value = 44413
$ws.Range("J44").Value = 700
$ws.Range("K44").Value = 26000
$ws.Range("L44").Value = 28000
$ws.Range("M44").Value = 27000
$ws.Range("P44").Value = 1080

# Row 45
$ws.Range("D45").Value = 44399
$ws.Range("J45").Value = 400

# Row 46
$ws.Range("D46").Value = 44349
$ws.Range("J46").Value = 600
